# Refactoring custom appearance demo into a separate form.
# This edits the "selects" XLSForm workbook:
#  - The old "inline" choose-one demo (i1/i2/i3, Choose one:) on the
#    survey sheet is replaced with a new "Which devices do you use?"
#    note + 4 inline select_one yes_no rows (desktop/laptop/smartphone/tablet).
#  - The content_provider_test row now has a condition of FALSE (hidden).
#  - The "queries" sheet's content_provider_test callback text changes
#    from "[context]" to "context".
#  - Column E width on survey sheet changes.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

# Insert two new rows right after the existing "i2"/"i3" rows (old rows 11-12),
# before the old "end screen" row (old row 13), so the appearance demo grows
# from 3 choices to 4 choices.
$survey.Rows.Item(13).Resize(2).Insert()

# Row 10: was "select_one yes_no / inline / i1 / Choose one:" -- now becomes
# a plain note introducing the new device question.
$survey.Range("B10").Value = "note"
$survey.Range("C10").Value = ""
$survey.Range("F10").Value = ""
$survey.Range("G10").Value = "Which devices do you use?"

# Rows 11-14: the four inline yes_no selects for the devices question.
$survey.Range("B11").Value = "select_one yes_no"
$survey.Range("C11").Value = "inline"
$survey.Range("F11").Value = "desktop"
$survey.Range("G11").Value = "Desktop computer"

$survey.Range("B12").Value = "select_one yes_no"
$survey.Range("C12").Value = "inline"
$survey.Range("F12").Value = "laptop"
$survey.Range("G12").Value = "Laptop computer"

$survey.Range("B13").Value = "select_one yes_no"
$survey.Range("C13").Value = "inline"
$survey.Range("F13").Value = "smartphone"
$survey.Range("G13").Value = "Smartphone"

$survey.Range("B14").Value = "select_one yes_no"
$survey.Range("C14").Value = "inline"
$survey.Range("F14").Value = "tablet"
$survey.Range("G14").Value = "Tablet"

# The content provider query demo (now on row 21 after the insert) is
# hidden by default via a FALSE relevance condition.
$survey.Range("E21").Value = $false

# Column E width changes slightly (narrower). The runtime quantizes
# ColumnWidth to sixths of a character, so 32.5 is the closest input that
# lands on the nearest representable stored width to the target 33.29.
$survey.Columns.Item(5).ColumnWidth = 32.5

# "queries" sheet: callback for content_provider_test changes from
# "[context]" to "context".
$queries = $wb.Worksheets.Item("queries")
$queries.Range("C6").Value = "context"

"done"
